$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.146.67'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').Value = '2.090.31'
$ws.Range('E3').Value = '  +8.95%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'251.16"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.78%  '
$ws.Range('D6').Value = "'0.656"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -4.45%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = "'50.55"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +6.31%  '
$ws.Range('D9').Value = "'60.81"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +4.50%  '
$ws.Range('D10').Value = "'0.372"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('D11').Value = "'0.0745"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -2.12%  '
$ws.Range('E12').Value = '  +5.23%  '
$ws.Range('D13').Value = "'15.16"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('E14').Value = '  +9.24%  '
$ws.Range('D15').Value = "'0.831"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +1.23%  '
$ws.Range('D16').Value = '2.092.27'
$ws.Range('E16').Value = '  +9.07%  '
$ws.Range('D17').Value = "'5.10"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').Value = '37.084.72'
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = "'72.32"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').Value = '0.0₃0822'
$ws.Range('E20').Value = '  -3.71%  '
$ws.Range('D21').Value = "'13.26"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.60%  '
$ws.Range('D22').Value = "'240.07"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').Value = "'5.21"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.12%  '
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').Value = "'2.47"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('D26').Value = "'169.29"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.05%  '
$ws.Range('D27').Value = "'9.27"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +5.41%  '
$ws.Range('D28').Value = "'21.12"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +12.97%  '
$ws.Range('D29').Value = "'2.00"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -4.59%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = "'0.122"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.72%  '
$ws.Range('B31').Value = 'Gas'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D31').Value = "'25.58"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +33.45%  '
$ws.Range('D32').Value = "'1.13"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +26.53%  '
$ws.Range('D33').Value = "'4.49"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('D34').Value = "'0.0606"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').Value = "'0.0930"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.33%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = "'2.28"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +17.05%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = "'1.84"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').Value = "'4.07"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.01%  '
$ws.Range('E40').Value = '  -9.28%  '
$ws.Range('D41').Value = "'17.63"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('D42').Value = "'0.0224"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('E43').Value = '  +5.24%  '
$ws.Range('D44').Value = "'97.56"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.91%  '
$ws.Range('D45').Value = "'2.81"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('D46').Value = "'0.0869"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.07%  '
$ws.Range('D47').Value = "'2.99"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +6.68%  '
$ws.Range('D48').Value = '1.308.43'
$ws.Range('E48').Value = '  -2.97%  '
$ws.Range('D49').Value = "'6.90"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +7.86%  '
$ws.Range('D50').Value = '2.279.29'
$ws.Range('E50').Value = '  +8.90%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = "'46.02"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.54%  '
